$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.484.21"
$ws.Cells.Item(2, 5).Value = "  +0.12%  "
$ws.Cells.Item(3, 4).Value = "1.913.36"
$ws.Cells.Item(3, 5).Value = "  -0.17%  "
$ws.Cells.Item(4, 4).Value = "'0.9989"
$ws.Cells.Item(4, 5).Value = "  -0.09%  "
$ws.Cells.Item(5, 4).Value = "'245.36"
$ws.Cells.Item(5, 5).Value = "  +1.49%  "
$ws.Cells.Item(6, 4).Value = "'0.9989"
$ws.Cells.Item(6, 5).Value = "  -0.12%  "
$ws.Cells.Item(7, 4).Value = "'0.4825"
$ws.Cells.Item(7, 5).Value = "  +2.62%  "
$ws.Cells.Item(8, 4).Value = "'0.2896"
$ws.Cells.Item(8, 5).Value = "  +1.25%  "
$ws.Cells.Item(9, 4).Value = "'0.06725"
$ws.Cells.Item(9, 5).Value = "  -1.07%  "
$ws.Cells.Item(10, 4).Value = "'110.79"
$ws.Cells.Item(10, 5).Value = "  +3.73%  "
$ws.Cells.Item(11, 4).Value = "'19.21"
$ws.Cells.Item(11, 5).Value = "  +4.33%  "
$ws.Cells.Item(12, 4).Value = "1.914.97"
$ws.Cells.Item(12, 5).Value = "  +0.05%  "
$ws.Cells.Item(13, 4).Value = "'0.07556"
$ws.Cells.Item(13, 5).Value = "  -2.03%  "
$ws.Cells.Item(14, 4).Value = "'5.277"
$ws.Cells.Item(14, 5).Value = "  +1.17%  "
$ws.Cells.Item(15, 4).Value = "'0.6718"
$ws.Cells.Item(15, 5).Value = "  +1.95%  "
$ws.Cells.Item(16, 4).Value = "'288.96"
$ws.Cells.Item(16, 5).Value = "  -0.35%  "
$ws.Cells.Item(17, 4).Value = "30.470.62"
$ws.Cells.Item(17, 5).Value = "  +0.05%  "
$ws.Cells.Item(18, 4).Value = "'0.000007600"
$ws.Cells.Item(18, 5).Value = "  -0.50%  "
$ws.Cells.Item(19, 4).Value = "'0.9988"
$ws.Cells.Item(20, 4).Value = "'12.88"
$ws.Cells.Item(20, 5).Value = "  -0.59%  "
$ws.Cells.Item(21, 4).Value = "2.159.45"
$ws.Cells.Item(21, 5).Value = "  +0.24%  "
$ws.Cells.Item(22, 4).Value = "'5.473"
$ws.Cells.Item(22, 5).Value = "  +5.02%  "
$ws.Cells.Item(23, 4).Value = "'0.9988"
$ws.Cells.Item(23, 5).Value = "  -0.17%  "
$ws.Cells.Item(24, 4).Value = "'6.404"
$ws.Cells.Item(24, 5).Value = "  +2.56%  "
$ws.Cells.Item(25, 4).Value = "'9.469"
$ws.Cells.Item(25, 5).Value = "  +1.51%  "
$ws.Cells.Item(26, 4).Value = "'164.11"
$ws.Cells.Item(26, 5).Value = "  -2.36%  "
$ws.Cells.Item(27, 4).Value = "'20.36"
$ws.Cells.Item(27, 5).Value = "  -5.38%  "
$ws.Cells.Item(28, 4).Value = "'2.131"
$ws.Cells.Item(28, 5).Value = "  +2.45%  "
$ws.Cells.Item(29, 4).Value = "'0.1058"
$ws.Cells.Item(29, 5).Value = "  -0.83%  "
$ws.Cells.Item(30, 4).Value = "'1.405"
$ws.Cells.Item(30, 5).Value = "  +2.67%  "
$ws.Cells.Item(31, 4).Value = "'4.184"
$ws.Cells.Item(31, 5).Value = "  +0.41%  "
$ws.Cells.Item(32, 4).Value = "'4.059"
$ws.Cells.Item(33, 4).Value = "'0.04994"
$ws.Cells.Item(33, 5).Value = "  -1.43%  "
$ws.Cells.Item(34, 4).Value = "'0.7299"
$ws.Cells.Item(34, 5).Value = "  -2.09%  "
$ws.Cells.Item(35, 4).Value = "'1.134"
$ws.Cells.Item(35, 5).Value = "  -1.62%  "
$ws.Cells.Item(36, 5).Value = "  -0.09%  "
# Row 37: HuobiToken -> VeChain
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(37, 4).Value = "'0.02047"
$ws.Cells.Item(37, 5).Value = "  -2.48%  "

# Row 38: VeChain -> HuobiToken
$ws.Cells.Item(38, 2).Value = "HuobiToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(38, 4).Value = "'2.718"
$ws.Cells.Item(38, 5).Value = "  -0.82%  "

$ws.Cells.Item(39, 4).Value = "'2.666"
$ws.Cells.Item(39, 5).Value = "  -0.16%  "
$ws.Cells.Item(40, 4).Value = "'110.79"
$ws.Cells.Item(40, 5).Value = "  +0.83%  "
$ws.Cells.Item(41, 4).Value = "'2.013"
$ws.Cells.Item(41, 5).Value = "  -2.00%  "
$ws.Cells.Item(42, 4).Value = "'0.4440"
$ws.Cells.Item(42, 5).Value = "  +4.30%  "
$ws.Cells.Item(43, 4).Value = "'0.8678"
$ws.Cells.Item(43, 5).Value = "  -0.25%  "
$ws.Cells.Item(44, 4).Value = "'5.842"
$ws.Cells.Item(44, 5).Value = "  -0.37%  "
$ws.Cells.Item(45, 5).Value = "  -0.12%  "
$ws.Cells.Item(46, 4).Value = "'68.22"
$ws.Cells.Item(46, 5).Value = "  +0.91%  "
$ws.Cells.Item(47, 4).Value = "'7.367"
$ws.Cells.Item(47, 5).Value = "  +2.85%  "
$ws.Cells.Item(48, 4).Value = "'49.00"
$ws.Cells.Item(48, 5).Value = "  -4.48%  "
$ws.Cells.Item(49, 4).Value = "'9.277"
$ws.Cells.Item(49, 5).Value = "  +0.01%  "
$ws.Cells.Item(50, 4).Value = "'0.1238"
$ws.Cells.Item(50, 5).Value = "  +1.94%  "
$ws.Cells.Item(51, 4).Value = "'34.87"
$ws.Cells.Item(51, 5).Value = "  -0.23%  "
